# Penambahan if else dan memperbaiki source code status
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix wording of status text values (shared-string text edits + value changes)
$ws.Range("F2").Value = "OK invalid id"
$ws.Range("E3").Value = "Akun Anda untuk sementara tidak dapat digunakan untuk berbelanja."
$ws.Range("F3").Value = "NO"
$ws.Range("F4").Value = "OK berhasil masuk profile"

# New row 5 with if/else style status addition
$ws.Range("E5").Value = "Akun Anda untuk sementara tidak dapat digunakan untuk berbelanja."
$ws.Range("F5").Value = "Oke"

# Update the hyperlinks so B3 uses rId1 (mailto) and B4 uses rId2 (invalid uri)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:hilmi.falih@yahoo.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "http://invalid.uri/") | Out-Null
$ws.Range("B3").Style = "Hyperlink"
$ws.Range("B4").Style = "Hyperlink"

# Update selection / view state
$ws.Range("D9").Select()
$excel.ActiveWindow.ScrollColumn = 2
